# Bereich 'Entgelt' in Datenbank eingefuegt.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tabelle1")
$ws2 = $wb.Worksheets.Item("Tabelle2")

# --- Tabelle2: add a new "Tarifbeschaeftigt?" ja/nein list in column I ---
$ws2.Range("I1").Value = "Tarifbeschaeftigt?"
$ws2.Range("I2").Value = "ja"
$ws2.Range("I3").Value = "nein"

# --- Tabelle1: fill in the already-existing "Tarif"/"Tarifbeschaeftigt?" rows ---
$ws1.Range("B33").Value = "nein"
$ws1.Range("B34").Value = "A5"

# Turn the old "aussertariflich beschaeftigt?" row into the new "Gewerkschaft" row,
# and restyle it like the "Tarif" row above it (fill style s=9).
$ws1.Range("A34:B34").Copy()
$ws1.Range("A35:B35").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws1.Range("A35").Value = "Gewerkschaft"
$ws1.Range("B35").Value = "verdi"

# Insert three new rows below it for the "Entgelt" figures and give them the
# same fill style (copy format from the row above, which is already s=9).
$ws1.Rows("36:38").Insert()
$ws1.Range("A35:D35").Copy()
$ws1.Range("A36:D38").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Range("A36").Value = "Grundgehalt"
$ws1.Range("B36").Value = 3523.76
$ws1.Range("A37").Value = "Weihnachtsgeld"
$ws1.Range("B37").Value = 1254.28
$ws1.Range("A38").Value = "Urlaubgeld"
$ws1.Range("B38").Value = 2568.75

# New dropdown validation on B33 pointing at the new Tabelle2 list.
$ws1.Range("B33").Validation.Add(3, 1, 1, "=Tabelle2!`$I`$2:`$I`$3")
$ws1.Range("B33").Validation.IgnoreBlank = $true
$ws1.Range("B33").Validation.InCellDropdown = $true

# --- view state bookkeeping (selection / scroll position) ---
$ws1.Application.ActiveWindow.ScrollRow = 22
$ws1.Range("B37").Select()
$ws2.Range("I5").Select()
$ws1.Select()
